# Fix casing of the "mdh_*" labels in column A to match the already-used
# "Mdh_*" capitalization style used elsewhere in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Mdh_2.5"
$ws.Range("A3").Value = "Mdh_1"
$ws.Range("A4").Value = "Mdh_.2"
$ws.Range("A5").Value = "Mdh_0"

# Reflect the author's last active selection in the sheet view.
$ws.Range("A6").Select()
